$d = $word.ActiveDocument

# The document originally starts with:
#   Para 1: "Supplement S2 Table"          (to be removed)
#   Para 2: "" (empty, bold formatting)    (to be removed)
#   Para 3: "S2 Table. Satellite Imagery." (stays; _GoBack bookmark currently
#            sits between the "2" and " Table." runs and must move to the
#            very start of this paragraph, right after the paragraph's pPr)
#
# Remove paragraph 2 (the empty paragraph) first. This keeps paragraph 1's
# content/runs completely untouched for now, so the bookmark relocation
# below can be anchored on a non-zero offset (position 0 triggers a bug in
# this runtime's Bookmarks.Add that unexpectedly expands the bookmark over
# the whole following paragraph instead of keeping it collapsed).
$d.Paragraphs.Item(2).Range.Delete()

# Relocate the "_GoBack" bookmark from its current spot (between the "2" and
# " Table." runs of what is now paragraph 2) to the very start of that
# paragraph's text.
$bookmark = $d.Bookmarks.Item("_GoBack")
$bookmark.Delete()

$targetParagraph = $d.Paragraphs.Item(2)
$startPos = $targetParagraph.Range.Start
$collapsedRange = $d.Range($startPos, $startPos)
$d.Bookmarks.Add("_GoBack", $collapsedRange)

# Finally remove paragraph 1 ("Supplement S2 Table"), merging the now-first
# paragraph (previously paragraph 2, with the relocated bookmark at its very
# start) up to become the document's first paragraph.
$d.Paragraphs.Item(1).Range.Delete()
